$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9460.916999999999
$ws.Range("I9").Value = 2029.5
$ws.Range("J9").Value = 13176.625
$ws.Range("K9").Value = 2029.5
$ws.Range("L9").Value = 13176.625
$ws.Range("M9").Value = -1860.5
$ws.Range("N9").Value = -13514.625

$ws.Range("H69").Value = 12380
$ws.Range("I69").Value = 12380
$ws.Range("K69").Value = 37140
$ws.Range("M69").Value = -36266

$ws.Range("H72").Value = 12380
$ws.Range("I72").Value = 12380
$ws.Range("K72").Value = 111420
$ws.Range("M72").Value = -107052

$ws.Range("H86").Value = 1851.4
$ws.Range("I86").Value = 1863.6666
$ws.Range("J86").Value = 1802.3334
$ws.Range("K86").Value = 1863.6666
$ws.Range("L86").Value = 1802.3334
$ws.Range("M86").Value = -740.6666
$ws.Range("N86").Value = -4048.3334

$ws.Range("H89").Value = 1851.4
$ws.Range("I89").Value = 1863.6666
$ws.Range("J89").Value = 1802.3334
$ws.Range("K89").Value = 9318.333000000001
$ws.Range("L89").Value = 9011.666999999999
$ws.Range("M89").Value = -3702.333000000001
$ws.Range("N89").Value = -20243.667

$ws.Range("H92").Value = 52342.965
$ws.Range("I92").Value = 28054.777
$ws.Range("K92").Value = 28054.777
$ws.Range("M92").Value = -26806.777

$ws.Range("H96").Value = 1338.5652
$ws.Range("I96").Value = 506.9
$ws.Range("J96").Value = 1978.3077
$ws.Range("K96").Value = 1520.7
$ws.Range("L96").Value = 5934.9231
$ws.Range("M96").Value = -147.6999999999998
$ws.Range("N96").Value = -8680.9231

$ws.Range("H97").Value = 200659.8
$ws.Range("J97").Value = 200659.8
$ws.Range("L97").Value = 601979.3999999999
$ws.Range("N97").Value = -602971.3999999999

$ws.Range("H99").Value = 278.07693
$ws.Range("I99").Value = 245.6
$ws.Range("K99").Value = 736.8
$ws.Range("M99").Value = 761.2

$ws.Range("H101").Value = 2127.5264
$ws.Range("I101").Value = 1709.4
$ws.Range("J101").Value = 2592.111
$ws.Range("K101").Value = 5128.200000000001
$ws.Range("L101").Value = 7776.333
$ws.Range("M101").Value = -3506.200000000001
$ws.Range("N101").Value = -11020.333

$ws.Range("H113").Value = 2961.5
$ws.Range("I113").Value = 2262
$ws.Range("J113").Value = 3427.8333
$ws.Range("K113").Value = 2262
$ws.Range("L113").Value = 3427.8333
$ws.Range("M113").Value = 992
$ws.Range("N113").Value = -9935.8333

$ws.Range("H132").Value = 2693.4773
$ws.Range("I132").Value = 2345.6191
$ws.Range("J132").Value = 9998.5
$ws.Range("K132").Value = 7036.8573
$ws.Range("L132").Value = 29995.5
$ws.Range("M132").Value = -4506.8573
$ws.Range("N132").Value = -35055.5

$ws.Range("H138").Value = 3607.25
$ws.Range("I138").Value = 1712
$ws.Range("K138").Value = 5136
$ws.Range("M138").Value = 4

$ws.Range("H141").Value = 397.0345
$ws.Range("I141").Value = 388.80356
$ws.Range("K141").Value = 1166.41068
$ws.Range("M141").Value = 4013.58932

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5857.4487
$ws.Range("I32").Value = 3780.6934
$ws.Range("K32").Value = 3780.6934
$ws.Range("M32").Value = -3493.6934

$ws.Range("H42").Value = 12000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 12000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 12000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -12972

$ws.Range("H74").Value = 348740.94
$ws.Range("I74").Value = 332036.12
$ws.Range("K74").Value = 332036.12
$ws.Range("M74").Value = -331162.12

$ws.Range("H77").Value = 348740.94
$ws.Range("I77").Value = 332036.12
$ws.Range("K77").Value = 1660180.6
$ws.Range("M77").Value = -1655812.6

$ws.Range("H102").Value = 324792.94
$ws.Range("I102").Value = 478264.53
$ws.Range("K102").Value = 478264.53
$ws.Range("M102").Value = -476642.53

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 24436626
$ws.Range("I86").Value = 37072710
$ws.Range("K86").Value = 37072710
$ws.Range("M86").Value = -37071587

$ws.Range("H89").Value = 24436626
$ws.Range("I89").Value = 37072710
$ws.Range("K89").Value = 185363550
$ws.Range("M89").Value = -185357934

$ws.Range("H94").Value = 30521.8
$ws.Range("I94").Value = 503.08334
$ws.Range("K94").Value = 503.08334
$ws.Range("M94").Value = -52.08334000000002

$ws.Range("H99").Value = 2147.75
$ws.Range("J99").Value = 2160.75
$ws.Range("L99").Value = 2160.75
$ws.Range("N99").Value = -5156.75

$ws.Range("H105").Value = 31253200
$ws.Range("I105").Value = 38464684
$ws.Range("J105").Value = 3436.6667
$ws.Range("K105").Value = 38464684
$ws.Range("L105").Value = 3436.6667
$ws.Range("M105").Value = -38462937
$ws.Range("N105").Value = -6930.6667

$ws.Range("H107").Value = 6057.4043
$ws.Range("I107").Value = 5404.8203
$ws.Range("J107").Value = 9238.75
$ws.Range("K107").Value = 5404.8203
$ws.Range("L107").Value = 9238.75
$ws.Range("M107").Value = -3484.8203
$ws.Range("N107").Value = -13078.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 176.53334
$ws.Range("I7").Value = 106.333336
$ws.Range("J7").Value = 281.83334
$ws.Range("K7").Value = 106.333336
$ws.Range("L7").Value = 281.83334
$ws.Range("M7").Value = 6.666663999999997
$ws.Range("N7").Value = -507.83334

$ws.Range("H11").Value = 797.5
$ws.Range("I11").Value = 397.5
$ws.Range("J11").Value = 997.5
$ws.Range("K11").Value = 397.5
$ws.Range("L11").Value = 997.5
$ws.Range("M11").Value = -257.5
$ws.Range("N11").Value = -1277.5

$ws.Range("H69").Value = 18635.4
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 18635.4
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H88").Value = 15988.4
$ws.Range("J88").Value = 14174.777
$ws.Range("L88").Value = 14174.777
$ws.Range("N88").Value = -14986.777

$ws.Range("H91").Value = 15988.4
$ws.Range("J91").Value = 14174.777
$ws.Range("L91").Value = 14174.777
$ws.Range("N91").Value = -16982.777

$ws.Range("H105").Value = 1720.0435
$ws.Range("I105").Value = 1328.5
$ws.Range("J105").Value = 2615
$ws.Range("K105").Value = 1328.5
$ws.Range("L105").Value = 2615
$ws.Range("M105").Value = 418.5
$ws.Range("N105").Value = -6109

$ws.Range("H112").Value = 67720.2
$ws.Range("J112").Value = 67720.2
$ws.Range("L112").Value = 67720.2
$ws.Range("N112").Value = -70674.2

$ws.Range("H122").Value = 1634.25
$ws.Range("I122").Value = 1345.6666
$ws.Range("K122").Value = 4036.9998
$ws.Range("M122").Value = -1586.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1494.9524
$ws.Range("I107").Value = 536.5
$ws.Range("J107").Value = 1595.8422
$ws.Range("K107").Value = 1609.5
$ws.Range("L107").Value = 4787.5266
$ws.Range("M107").Value = 310.5
$ws.Range("N107").Value = -8627.526600000001

$ws.Range("H140").Value = 2127.2195
$ws.Range("J140").Value = 2601.182
$ws.Range("L140").Value = 7803.545999999999
$ws.Range("N140").Value = -18163.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 30412.62
$ws.Range("J102").Value = 14419.583
$ws.Range("L102").Value = 14419.583
$ws.Range("N102").Value = -17663.583

$ws.Range("H122").Value = 259070.86
$ws.Range("I122").Value = 450999
$ws.Range("K122").Value = 1352997
$ws.Range("M122").Value = -1350547

$ws.Range("H126").Value = 51247.26
$ws.Range("I126").Value = 63886.332
$ws.Range("J126").Value = 3850.75
$ws.Range("K126").Value = 191658.996
$ws.Range("L126").Value = 11552.25
$ws.Range("M126").Value = -189188.996
$ws.Range("N126").Value = -16492.25

$ws.Range("H132").Value = 5516.8125
$ws.Range("I132").Value = 5644.6665
$ws.Range("J132").Value = 3599
$ws.Range("K132").Value = 16933.9995
$ws.Range("L132").Value = 10797
$ws.Range("M132").Value = -14403.9995
$ws.Range("N132").Value = -15857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1356.2174
$ws.Range("I93").Value = 1190
$ws.Range("J93").Value = 2145.75
$ws.Range("K93").Value = 1190
$ws.Range("L93").Value = 2145.75
$ws.Range("M93").Value = 58
$ws.Range("N93").Value = -4641.75

$ws.Range("H122").Value = 4308.1
$ws.Range("I122").Value = 4120.222
$ws.Range("K122").Value = 12360.666
$ws.Range("M122").Value = -9910.665999999999

$ws.Range("H136").Value = 51263.74
$ws.Range("I136").Value = 3319.9092
$ws.Range("K136").Value = 9959.7276
$ws.Range("M136").Value = -7409.7276

$ws.Range("H141").Value = 125000
$ws.Range("J141").Value = 125000
$ws.Range("L141").Value = 125000
$ws.Range("N141").Value = -135360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 30239.8
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 30239.8
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 30239.8
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -30995.8

$ws.Range("H113").Value = 917.5185
$ws.Range("I113").Value = 868.1923
$ws.Range("K113").Value = 2604.5769
$ws.Range("M113").Value = -434.5769
